$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 16) mirroring the existing data rows, reusing the
# "HexGrid-60degTilt5degRes" label (same as row 15) with row index 14 in column A.
# Copy the formatting of A15 (bold, centered, bordered) onto A16.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.995896977005419
$ws.Range("D16").Value = 0.9884528854690178
$ws.Range("E16").Value = 0.9946551405311592
$ws.Range("F16").Value = 0.995896977005419
$ws.Range("G16").Value = 0.9812925782960419
$ws.Range("H16").Value = 0.9923753196759012
$ws.Range("I16").Value = 0.9906682252443847
$ws.Range("J16").Value = 0.9884528854690178
$ws.Range("K16").Value = 0.9915540130000885
$ws.Range("L16").Value = 0.9937254950027536
$ws.Range("M16").Value = 0.9905568543703205
